# Weekly price update: insert a new "Pepino ensalada" price record for
# Vega Central Mapocho de Santiago at row 175, pushing the existing
# rows 175-249 down to 176-250 (dimension grows from A1:R249 to A1:R250).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 175..249 down one row, creating a blank row 175.
$ws.Rows.Item(175).Insert()

# Populate the newly inserted row 175 with the new weekly record.
$ws.Cells.Item(175, 1).Value  = 9
$ws.Cells.Item(175, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(175, 3).Value  = "Metropolitana"
$ws.Cells.Item(175, 4).Value  = 44627
$ws.Cells.Item(175, 5).Value  = 13
$ws.Cells.Item(175, 6).Value  = 100112043
$ws.Cells.Item(175, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(175, 8).Value  = "Sin especificar"
$ws.Cells.Item(175, 9).Value  = "Primera"
$ws.Cells.Item(175, 10).Value = 52
$ws.Cells.Item(175, 11).Value = 17000
$ws.Cells.Item(175, 12).Value = 18000
$ws.Cells.Item(175, 13).Value = 17500
$ws.Cells.Item(175, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(175, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(175, 16).Value = 292
$ws.Cells.Item(175, 17).Value = 60
$ws.Cells.Item(175, 18).Value = "Hortaliza"
